$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_7_3_22"
$ws.Range("B2").Value = 0.6527417204065737
$ws.Range("C2").Value = -1.076986097008408
$ws.Range("D2").Value = 0.2831971299617413
$ws.Range("E2").Value = -0.4340917799518298
$ws.Range("F2").Value = 0.3843125998973846
$ws.Range("G2").Value = 2.002768516540527
$ws.Range("H2").Value = 0.3571963012218475
$ws.Range("I2").Value = 1.228379487991333

$ws.Range("A3").Value = "model_7_3_23"
$ws.Range("B3").Value = 0.6557679162826306
$ws.Range("C3").Value = -1.059926389773422
$ws.Range("D3").Value = 0.2807720278777682
$ws.Range("E3").Value = -0.424593036065892
$ws.Range("F3").Value = 0.3809634745121002
$ws.Range("G3").Value = 1.986318349838257
$ws.Range("H3").Value = 0.3584047853946686
$ws.Range("I3").Value = 1.220243334770203

$ws.Range("A4").Value = "model_7_3_24"
$ws.Range("B4").Value = 0.6581054898367903
$ws.Range("C4").Value = -1.042470777745613
$ws.Range("D4").Value = 0.2702141518763468
$ws.Range("E4").Value = -0.4170822945868511
$ws.Range("F4").Value = 0.3783764839172363
$ws.Range("G4").Value = 1.969486713409424
$ws.Range("H4").Value = 0.3636659979820251
$ws.Range("I4").Value = 1.213809967041016

$ws.Range("A5").Value = "model_7_3_21"
$ws.Range("B5").Value = 0.6616847229553402
$ws.Range("C5").Value = -1.025443730758705
$ws.Range("D5").Value = 0.3112407722137548
$ws.Range("E5").Value = -0.3956982949722045
$ws.Range("F5").Value = 0.3744153380393982
$ws.Range("G5").Value = 1.953068017959595
$ws.Range("H5").Value = 0.3432216644287109
$ws.Range("I5").Value = 1.195493340492249

$ws.Range("A6").Value = "model_7_3_20"
$ws.Range("B6").Value = 0.6622821195235455
$ws.Range("C6").Value = -1.015258006641595
$ws.Range("D6").Value = 0.3036606396937693
$ws.Range("E6").Value = -0.3917050613535327
$ws.Range("F6").Value = 0.3737541735172272
$ws.Range("G6").Value = 1.943246245384216
$ws.Range("H6").Value = 0.3469989597797394
$ws.Range("I6").Value = 1.192072868347168

$ws.Range("A7").Value = "model_7_3_0"
$ws.Range("B7").Value = 0.6656520100818517
$ws.Range("C7").Value = 0.4328243385309309
$ws.Range("D7").Value = 0.1968968018885028
$ws.Range("E7").Value = 0.4421024141572517
$ws.Range("F7").Value = 0.3700246810913086
$ws.Range("G7").Value = 0.5469086170196533
$ws.Range("H7").Value = 0.4002013802528381
$ws.Range("I7").Value = 0.477870374917984

$ws.Range("A8").Value = "model_7_3_18"
$ws.Range("B8").Value = 0.6697019019140396
$ws.Range("C8").Value = -0.9574554003626063
$ws.Range("D8").Value = 0.3108935369669598
$ws.Range("E8").Value = -0.3552710248468809
$ws.Range("F8").Value = 0.3655426502227783
$ws.Range("G8").Value = 1.887508988380432
$ws.Range("H8").Value = 0.3433946967124939
$ws.Range("I8").Value = 1.160865187644958

$ws.Range("A9").Value = "model_7_3_19"
$ws.Range("B9").Value = 0.6723874611275902
$ws.Range("C9").Value = -0.948376761427969
$ws.Range("D9").Value = 0.3194469156731838
$ws.Range("E9").Value = -0.3475185608506703
$ws.Range("F9").Value = 0.3625705540180206
$ws.Range("G9").Value = 1.87875497341156
$ws.Range("H9").Value = 0.3391323685646057
$ws.Range("I9").Value = 1.154224872589111

$ws.Range("A10").Value = "model_7_3_1"
$ws.Range("B10").Value = 0.6732547245156322
$ws.Range("C10").Value = 0.3861529080667694
$ws.Range("D10").Value = 0.1187722819988261
$ws.Range("E10").Value = 0.3928961228641034
$ws.Range("F10").Value = 0.3616107106208801
$ws.Range("G10").Value = 0.5919123291969299
$ws.Range("H10").Value = 0.4391322731971741
$ws.Range("I10").Value = 0.5200182795524597

$ws.Range("A11").Value = "model_7_3_4"
$ws.Range("B11").Value = 0.6738025962167786
$ws.Range("C11").Value = -0.4708913753111801
$ws.Range("D11").Value = -0.03590516389094134
$ws.Range("E11").Value = -0.160232878260385
$ws.Range("F11").Value = 0.3610044121742249
$ws.Range("G11").Value = 1.418331623077393
$ws.Range("H11").Value = 0.5162109732627869
$ws.Range("I11").Value = 0.9938041567802429

$ws.Range("A12").Value = "model_7_3_17"
$ws.Range("B12").Value = 0.6765076320527268
$ws.Range("C12").Value = -0.9062226958472819
$ws.Range("D12").Value = 0.3113211583975165
$ws.Range("E12").Value = -0.3246199529532525
$ws.Range("F12").Value = 0.3580107092857361
$ws.Range("G12").Value = 1.838106989860535
$ws.Range("H12").Value = 0.3431815803050995
$ws.Range("I12").Value = 1.134610772132874

$ws.Range("A13").Value = "model_7_3_15"
$ws.Range("B13").Value = 0.6831351593780601
$ws.Range("C13").Value = -0.8158529917709556
$ws.Range("D13").Value = 0.2632816383256307
$ws.Range("E13").Value = -0.2839174171709842
$ws.Range("F13").Value = 0.3506760001182556
$ws.Range("G13").Value = 1.750966548919678
$ws.Range("H13").Value = 0.3671205639839172
$ws.Range("I13").Value = 1.099746704101562

$ws.Range("A14").Value = "model_7_3_3"
$ws.Range("B14").Value = 0.6832342214041476
$ws.Range("C14").Value = -0.2789527786092287
$ws.Range("D14").Value = -0.06217504647952543
$ws.Range("E14").Value = -0.05303254264701351
$ws.Range("F14").Value = 0.3505663871765137
$ws.Range("G14").Value = 1.233251690864563
$ws.Range("H14").Value = 0.5293017625808716
$ws.Range("I14").Value = 0.9019810557365417

$ws.Range("A15").Value = "model_7_3_16"
$ws.Range("B15").Value = 0.6877769576055393
$ws.Range("C15").Value = -0.78636728889883
$ws.Range("D15").Value = 0.2607254164085511
$ws.Range("E15").Value = -0.2670402084037933
$ws.Range("F15").Value = 0.3455388844013214
$ws.Range("G15").Value = 1.722534537315369
$ws.Range("H15").Value = 0.3683943748474121
$ws.Range("I15").Value = 1.085290551185608

$ws.Range("A16").Value = "model_7_3_2"
$ws.Range("B16").Value = 0.7070344925899879
$ws.Range("C16").Value = -0.1539898299829319
$ws.Range("D16").Value = 0.05940182219026624
$ws.Range("E16").Value = 0.05472783179540841
$ws.Range("F16").Value = 0.3242264688014984
$ws.Range("G16").Value = 1.112753987312317
$ws.Range("H16").Value = 0.4687176942825317
$ws.Range("I16").Value = 0.809678316116333

$ws.Range("A17").Value = "model_7_3_11"
$ws.Range("B17").Value = 0.7223410045836549
$ws.Range("C17").Value = -0.5244593641870401
$ws.Range("D17").Value = 0.3622080320478247
$ws.Range("E17").Value = -0.08316649438548107
$ws.Range("F17").Value = 0.3072867095470428
$ws.Range("G17").Value = 1.469985485076904
$ws.Range("H17").Value = 0.3178237080574036
$ws.Range("I17").Value = 0.9277924299240112

$ws.Range("A18").Value = "model_7_3_10"
$ws.Range("B18").Value = 0.7249347593765983
$ws.Range("C18").Value = -0.5012871446098954
$ws.Range("D18").Value = 0.3664222266390378
$ws.Range("E18").Value = -0.06820142639812365
$ws.Range("F18").Value = 0.3044161796569824
$ws.Range("G18").Value = 1.447641253471375
$ws.Range("H18").Value = 0.3157236874103546
$ws.Range("I18").Value = 0.9149740934371948

$ws.Range("A19").Value = "model_7_3_9"
$ws.Range("B19").Value = 0.7258458036428029
$ws.Range("C19").Value = -0.4842655084454499
$ws.Range("D19").Value = 0.3595989038938505
$ws.Range("E19").Value = -0.05992451687219313
$ws.Range("F19").Value = 0.3034079372882843
$ws.Range("G19").Value = 1.431227803230286
$ws.Range("H19").Value = 0.3191238641738892
$ws.Range("I19").Value = 0.907884418964386

$ws.Range("A20").Value = "model_7_3_14"
$ws.Range("B20").Value = 0.7347336407115015
$ws.Range("C20").Value = -0.5178306739768346
$ws.Range("D20").Value = 0.4525841892862353
$ws.Range("E20").Value = -0.05447235213876733
$ws.Range("F20").Value = 0.2935717105865479
$ws.Range("G20").Value = 1.463593602180481
$ws.Range("H20").Value = 0.2727875709533691
$ws.Range("I20").Value = 0.9032142758369446

$ws.Range("A21").Value = "model_7_3_8"
$ws.Range("B21").Value = 0.7370748460540743
$ws.Range("C21").Value = -0.3702127409088514
$ws.Range("D21").Value = 0.3546238795574079
$ws.Range("E21").Value = 0.00668580879056424
$ws.Range("F21").Value = 0.2909806668758392
$ws.Range("G21").Value = 1.321250438690186
$ws.Range("H21").Value = 0.321603000164032
$ws.Range("I21").Value = 0.8508289456367493

$ws.Range("A22").Value = "model_7_3_7"
$ws.Range("B22").Value = 0.7386364994529717
$ws.Range("C22").Value = -0.3528109894226359
$ws.Range("D22").Value = 0.3585588922394196
$ws.Range("E22").Value = 0.01813520385558787
$ws.Range("F22").Value = 0.289252370595932
$ws.Range("G22").Value = 1.304470658302307
$ws.Range("H22").Value = 0.3196421265602112
$ws.Range("I22").Value = 0.8410219550132751

$ws.Range("A23").Value = "model_7_3_5"
$ws.Range("B23").Value = 0.7411201802507892
$ws.Range("C23").Value = -0.3186781200121902
$ws.Range("D23").Value = 0.3605813847541544
$ws.Range("E23").Value = 0.03903198611842817
$ws.Range("F23").Value = 0.2865037024021149
$ws.Range("G23").Value = 1.271557331085205
$ws.Range("H23").Value = 0.3186342716217041
$ws.Range("I23").Value = 0.8231225609779358

$ws.Range("A24").Value = "model_7_3_6"
$ws.Range("B24").Value = 0.7429013378154702
$ws.Range("C24").Value = -0.3239783361199589
$ws.Range("D24").Value = 0.3790390922539143
$ws.Range("E24").Value = 0.04092549056512562
$ws.Range("F24").Value = 0.2845324873924255
$ws.Range("G24").Value = 1.276668190956116
$ws.Range("H24").Value = 0.3094364702701569
$ws.Range("I24").Value = 0.8215007185935974

$ws.Range("A25").Value = "model_7_3_13"
$ws.Range("B25").Value = 0.7440468795814711
$ws.Range("C25").Value = -0.4696435426386494
$ws.Range("D25").Value = 0.4992462333083746
$ws.Range("E25").Value = -0.01297832055689629
$ws.Range("F25").Value = 0.2832646667957306
$ws.Range("G25").Value = 1.417128324508667
$ws.Range("H25").Value = 0.2495349943637848
$ws.Range("I25").Value = 0.8676723837852478

$ws.Range("A26").Value = "model_7_3_12"
$ws.Range("B26").Value = 0.746085185891634
$ws.Range("C26").Value = -0.4537914509202188
$ws.Range("D26").Value = 0.5075222197061466
$ws.Range("E26").Value = -0.001267400542537933
$ws.Range("F26").Value = 0.2810088992118835
$ws.Range("G26").Value = 1.401842713356018
$ws.Range("H26").Value = 0.2454109191894531
$ws.Range("I26").Value = 0.857641339302063
